$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Row 1: drop the extra "value" labels in C1:F1 (only A1/B1 remain)
$ws.Range("C1:F1").ClearContents()

# Insert a new row above the old row 9 ("estimate_params") to make room for
# the new "L_curve" parameter row, pushing estimate_params..Deletion down by one.
$ws.Rows("9:9").Insert()

# Row 8 is now relabeled "production_function" (was "Model"); keep its value "Sigmoid".
$ws.Range("A8").Value = "production_function"
$ws.Range("A8").Style = $ws.Range("A1").Style

# Fill the newly-inserted row 9 with the L_curve parameter.
$ws.Range("A9").Value = "L_curve"
$ws.Range("A9").Style = $ws.Range("A1").Style
$ws.Range("B9").Value = 1
$ws.Range("B9").Style = $ws.Range("B2").Style

# The old "Deletion" row (now shifted to row 17) is no longer used - remove it entirely.
$ws.Rows("17:17").Delete()

# Update the sheet's remembered selection to match the authored file.
$ws.Range("B10").Select()
